$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1 (rows 235-239): rotate the match-data columns (B:AC) so the
# row that used to sit at 235 ends up at 239, with 236-239 each shifting up
# by one row. Column A (the sequential id) stays put on every row.
$buffer = $ws.Range("B235:AC235").Value2

$ws.Range("B235:AC235").Value2 = $ws.Range("B236:AC236").Value2
$ws.Range("B236:AC236").Value2 = $ws.Range("B237:AC237").Value2
$ws.Range("B237:AC237").Value2 = $ws.Range("B238:AC238").Value2
$ws.Range("B238:AC238").Value2 = $ws.Range("B239:AC239").Value2
$ws.Range("B239:AC239").Value2 = $buffer

# --- Block 2 (rows 267-274): the two earliest fixtures (old rows 267 and
# 268) are dropped. The remaining fixtures' data (B:AC) shifts up by two
# rows (269->267 ... 274->272), again leaving column A untouched, then the
# now-duplicate trailing rows 273:274 are removed outright so the sheet
# shrinks from 274 to 272 rows.
$ws.Range("B267:AC267").Value2 = $ws.Range("B269:AC269").Value2
$ws.Range("B268:AC268").Value2 = $ws.Range("B270:AC270").Value2
$ws.Range("B269:AC269").Value2 = $ws.Range("B271:AC271").Value2
$ws.Range("B270:AC270").Value2 = $ws.Range("B272:AC272").Value2
$ws.Range("B271:AC271").Value2 = $ws.Range("B273:AC273").Value2
$ws.Range("B272:AC272").Value2 = $ws.Range("B274:AC274").Value2

$ws.Range("A273:A274").EntireRow.Delete()
